# Set CRS for points_all_cleaned
# The cleaning step re-ran with an updated CRS transform, which dropped the
# first three TumTum "Wetland-Polygon" rows from the points table. Removing
# those rows shifts the remaining TumTum rows up by three and shrinks the
# used range accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A280:E282").EntireRow.Delete()
